# Add "4 pin SIPP socket" part to the Water Leak Alarm parts list.
#
# This inserts a new row (row 24) on Sheet1, pushing every row below it
# down by one, fills in the data for the newly added part, repairs a
# couple of relative-reference formulas that the row-insert leaves
# slightly wrong, restores the two hyperlinks that used to live at
# F24/F25 to their new homes at F25/F26 (together with their original,
# non-"Hyperlink"-styled formatting), and finally updates the sheet's
# active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

function Get-HyperlinkAddressAt($ws, $addrTarget) {
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq $addrTarget) {
            return $hl.Address
        }
    }
    return $null
}
function Remove-HyperlinkAt($ws, $addrTarget) {
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq $addrTarget) {
            $hl.Delete()
            return
        }
    }
}

# --- 1. Insert a new row at 24 (shifts old rows 24.. down to 25..) ---
$ws.Rows.Item(24).Insert()

# --- 2. The hyperlinks that used to sit at F24/F25 now belong to the
#        rows that were pushed down to F25/F26; remember & remove them
#        so we can re-create them in the right place below. ---
$oldF24Link = Get-HyperlinkAddressAt $ws "`$F`$24"
$oldF25Link = Get-HyperlinkAddressAt $ws "`$F`$25"
Remove-HyperlinkAt $ws "`$F`$24"
Remove-HyperlinkAt $ws "`$F`$25"

# --- 3. Populate the new row 24: "4 pin SIPP socket" from Jameco,
#        used to mount the DHT11 sensor ---
$ws.Cells.Item(24, 1).Formula = "=A23+1"
$ws.Cells.Item(24, 2).Value = 1
$ws.Cells.Item(24, 3).Value = "Jameco"
$ws.Cells.Item(24, 4).Value = 164822
$ws.Cells.Item(24, 5).Value = "4 pin SIPP socket"
$ws.Cells.Item(24, 6).Value = "http://www.jameco.com/z/6100-1-4-Socket-SIPP-1x4-Pin-Machine-Tool-Pins-Soldertail-Female_164822.html"
$ws.Cells.Item(24, 7).Value = 0.59
$ws.Cells.Item(24, 8).Formula = "=B24*G24"
$ws.Cells.Item(24, 9).Value = "Best to mount the DHT11 in a socket."
$ws.Rows.Item(24).RowHeight = 53.25

# --- 4. Re-create the displaced hyperlinks one row further down ---
if ($oldF24Link) { $ws.Hyperlinks.Add($ws.Range("F25"), $oldF24Link) | Out-Null }
if ($oldF25Link) { $ws.Hyperlinks.Add($ws.Range("F26"), $oldF25Link) | Out-Null }

# Adding a hyperlink applies the blue/underlined "Hyperlink" cell style;
# the original sheet instead keeps the plain wrapped-text look, so copy
# the formatting back from a cell that was never touched.
$ws.Range("F23").Copy()
$ws.Range("F25").PasteSpecial(-4122)
$ws.Range("F26").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- 5. Repair the running "Item #" counter / extended-price formulas on
#        the rows immediately surrounding the insertion point (the
#        automatic reference adjustment for the row right below the new
#        one needs a nudge) ---
$ws.Cells.Item(25, 1).Formula = "=A24+1"
$ws.Cells.Item(26, 1).Formula = "=A25+1"
$ws.Cells.Item(26, 8).Formula = "=B26*G26"

$excel.Calculate()

# --- 6. Update the sheet's active selection to reflect where the editor
#        left off after adding the new part ---
$ws.Range("H28").Select()
